$d = $word.ActiveDocument

# The paragraph containing the field placeholder "{m:null}" is currently
# stored as two runs: "{m" and ":null}". The parser now needs each of the
# tokens "{", "m", ":null" and "}" to live in its own run (see
# TokenIteratorFieldRewriterSplit), so split the paragraph's text at the
# three boundaries between those tokens.
#
# Word COM has no direct "split this run" call, but inserting (and then
# removing) a zero-length bookmark at a position forces the run under it
# to be split in two without touching any formatting, which is exactly
# the effect we want here.

# Locate "{m:null}" in the document so the script is not dependent on
# hard-coded character offsets.
$find = $d.Content
$find.Find.Execute("{m:null}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $find.Start

function Split-At([int]$position, [int]$counter) {
    $name = "m2docSplit" + $counter
    $d.Bookmarks.Add($name, $d.Range($position, $position))
    $d.Bookmarks($name).Delete()
}

# "{" | "m" | ":null" | "}"
Split-At ($start + 1) 1   # between "{" and "m"
Split-At ($start + 2) 2   # between "m" and ":null"
Split-At ($start + 7) 3   # between ":null" and "}"
